$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.0498220640569395
$wsSummary.Range("C2").Value = 0.0498220640569395
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.09491525423728814
$wsSummary.Range("F2").Value = 0.2077151335311573
$wsSummary.Range("G2").Value = 0.5768621236133122
$wsSummary.Range("H2").Value = 0.7874197431781701
$wsSummary.Range("I2").Value = 28
$wsSummary.Range("J2").Value = 534
$wsSummary.Range("K2").Value = 0
$wsSummary.Range("L2").Value = 0

# --- Classification Report sheet ---
$wsClassification = $wb.Worksheets.Item("Classification Report")
$wsClassification.Range("B2").Value = 0
$wsClassification.Range("C2").Value = 0
$wsClassification.Range("D2").Value = 0

$wsClassification.Range("B3").Value = 0.0498220640569395
$wsClassification.Range("C3").Value = 1
$wsClassification.Range("D3").Value = 0.09491525423728814

$wsClassification.Range("B4").Value = 0.0498220640569395
$wsClassification.Range("C4").Value = 0.0498220640569395
$wsClassification.Range("D4").Value = 0.0498220640569395
$wsClassification.Range("E4").Value = 0.0498220640569395

$wsClassification.Range("B5").Value = 0.02491103202846975
$wsClassification.Range("C5").Value = 0.5
$wsClassification.Range("D5").Value = 0.04745762711864407

$wsClassification.Range("B6").Value = 0.002482238066893783
$wsClassification.Range("C6").Value = 0.0498220640569395
$wsClassification.Range("D6").Value = 0.004728873876590867

# --- Confusion Matrix sheet ---
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")
$wsConfusion.Range("B2").Value = 0
$wsConfusion.Range("C2").Value = 534

$wsConfusion.Range("B3").Value = 0
$wsConfusion.Range("C3").Value = 28
